$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.039.37"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "3.400.67"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.19%  "
$ws.Range("D7").Value = "3.400.74"
$ws.Range("E7").Value = "  -0.96%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.52"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("E12").Value = "  +2.14%  "
$ws.Range("D13").Value = "3.981.59"
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.125"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000171"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.62%  "
$ws.Range("D17").Value = "3.402.30"
$ws.Range("E17").Value = "  -1.02%  "
$ws.Range("D18").Value = "61.118.33"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("E20").Value = "  -1.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "387.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.560"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.75%  "
$ws.Range("E25").Value = "  +0.51%  "
$ws.Range("E26").Value = "  -3.41%  "
$ws.Range("D27").Value = "3.539.66"
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("E28").Value = "  +0.43%  "
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.43"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.73%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.82%  "
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("E33").Value = "  -3.61%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.65"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.52%  "
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "166.94"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("D38").Value = "3.432.30"
$ws.Range("E38").Value = "  -0.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.91%  "
$ws.Range("E40").Value = "  -3.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "28.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0779"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("E43").Value = "  -2.49%  "
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.84%  "
$ws.Range("E47").Value = "  -3.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.99%  "
$ws.Range("D49").Value = "2.502.63"
$ws.Range("E49").Value = "  -3.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.14%  "
